$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep their text formatting (values look numeric),
# so Excel does not auto-convert them to Number types.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.528.60"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
$ws.Range("D3").Value = "1.846.17"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "263.69"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").Value = "0.5219"
$ws.Range("E7").Value = "  +1.35%  "

# Row 8
$ws.Range("D8").Value = "0.3235"
$ws.Range("E8").Value = "  -0.53%  "

# Row 9
$ws.Range("D9").Value = "0.06801"
$ws.Range("E9").Value = "  +0.52%  "

# Row 10
$ws.Range("D10").Value = "18.76"
$ws.Range("E10").Value = "  -1.22%  "

# Row 11
$ws.Range("D11").Value = "0.7771"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12
$ws.Range("D12").Value = "0.07788"
$ws.Range("E12").Value = "  +1.04%  "

# Row 13
$ws.Range("D13").Value = "1.844.70"
$ws.Range("E13").Value = "  -0.40%  "

# Row 14
$ws.Range("D14").Value = "88.41"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").Value = "5.017"

# Row 16
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("D17").Value = "13.97"
$ws.Range("E17").Value = "  -1.05%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007960"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "0.9995"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("D20").Value = "26.571.32"
$ws.Range("E20").Value = "  +0.22%  "

# Row 21
$ws.Range("D21").Value = "2.084.53"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22
$ws.Range("E22").Value = "  +1.77%  "

# Row 23
$ws.Range("D23").Value = "9.444"
$ws.Range("E23").Value = "  -1.07%  "

# Row 24
$ws.Range("D24").Value = "5.997"
$ws.Range("E24").Value = "  +0.81%  "

# Row 25
$ws.Range("D25").Value = "143.00"
$ws.Range("E25").Value = "  -1.10%  "

# Row 26
$ws.Range("D26").Value = "2.163"
$ws.Range("E26").Value = "  -8.11%  "

# Row 27
$ws.Range("D27").Value = "1.676"
$ws.Range("E27").Value = "  +1.22%  "

# Row 28
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").Value = "111.86"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30
$ws.Range("D30").Value = "4.183"
$ws.Range("E30").Value = "  -0.51%  "

# Row 31
$ws.Range("D31").Value = "0.08741"

# Row 32
$ws.Range("D32").Value = "4.110"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("D33").Value = "0.04831"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").Value = "0.7216"
$ws.Range("E34").Value = "  +4.76%  "

# Row 35
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
$ws.Range("D36").Value = "2.861"
$ws.Range("E36").Value = "  +0.42%  "

# Row 37
$ws.Range("D37").Value = "3.100"
$ws.Range("E37").Value = "  -0.60%  "

# Row 38
$ws.Range("D38").Value = "0.01792"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("D39").Value = "2.208"
$ws.Range("E39").Value = "  -0.62%  "

# Row 40
$ws.Range("D40").Value = "0.4844"

# Row 41
$ws.Range("D41").Value = "111.16"
$ws.Range("E41").Value = "  -1.78%  "

# Row 42
$ws.Range("D42").Value = "0.8916"
$ws.Range("E42").Value = "  -1.11%  "

# Row 43
$ws.Range("D43").Value = "6.034"
$ws.Range("E43").Value = "  -1.64%  "

# Row 44
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").Value = "7.617"
$ws.Range("E45").Value = "  -2.17%  "

# Row 46
$ws.Range("D46").Value = "0.4201"
$ws.Range("E46").Value = "  -0.98%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.083"
$ws.Range("E47").Value = "  -0.32%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05890"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49
$ws.Range("D49").Value = "0.1237"
$ws.Range("E49").Value = "  -2.40%  "

# Row 50
$ws.Range("D50").Value = "34.97"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("D51").Value = "0.8874"
$ws.Range("E51").Value = "  +3.99%  "
